$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 27 data
$ws.Range("A27").Value = "Payment"
$ws.Range("B27").Value = "bank defined beneficiaries"
$ws.Range("C27").Value = "ReferenceData.yaml"
$ws.Range("D27").Value = "GET"
$ws.Range("E27").Value = "/bankdefinedbeneficiaries"
$ws.Range("F27").Value = "In response, BDFResponse has property as BDFList. BDFList should have data type as Array."

# Adjust row 19 height (234 -> 187.2)
$ws.Rows("19").RowHeight = 187.2

# Update selection to new active cell
$ws.Range("A27:F27").Select()
